# Updated cryptos list on Mon May 22 12:29:11 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to set a text cell that Excel would otherwise auto-convert to a
# number (plain decimals like "312.29"). Forcing NumberFormat to "@" (Text)
# before the assignment keeps the literal string; resetting the Style back
# to "Normal" afterwards drops the now-unneeded cell style so the cell
# keeps the workbook's default (unstyled) look, matching the source data.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

function Set-Row($row, $b, $c, $d, $e, $dIsNumericLooking) {
    if ($null -ne $b) { $ws.Range("B$row").Value = $b }
    if ($null -ne $c) { $ws.Range("C$row").Value = $c }
    if ($null -ne $d) {
        if ($dIsNumericLooking) {
            Set-TextValue "D$row" $d
        } else {
            $ws.Range("D$row").Value = $d
        }
    }
    if ($null -ne $e) { $ws.Range("E$row").Value = $e }
}

Set-Row 2  $null $null "27.035.49" "  -0.13%  " $false
Set-Row 3  $null $null "1.829.88"  "  +0.21%  " $false
Set-Row 4  $null $null $null       "  -0.12%  " $false
Set-Row 5  $null $null "312.29"    "  -0.13%  " $true
Set-Row 6  $null $null "1.008"     "  -0.10%  " $true
Set-Row 7  $null $null "0.4659"    "  -0.71%  " $true
Set-Row 8  $null $null "0.3715"    "  +1.81%  " $true
Set-Row 9  $null $null "0.07386"   "  +0.08%  " $true
Set-Row 10 $null $null "0.8758"    "  -0.52%  " $true
Set-Row 11 $null $null "20.02"     "  -1.36%  " $true
Set-Row 12 $null $null "0.07806"   "  +6.74%  " $true
Set-Row 13 $null $null "1.846.91"  "  -4.56%  " $false

# Row 14 and 15 swap: Polkadot/Chainlink -> Chainlink/Polkadot
Set-Row 14 "Chainlink" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link" "6.593" "  +1.08%  " $true
Set-Row 15 "Polkadot" "https://coinranking.com/coin/25W7FG7om+polkadot-dot" "5.358" "  -0.14%  " $true

Set-Row 16 $null $null "91.92"         "  -1.44%  " $true
Set-Row 17 $null $null $null           "  +0.23%  " $false
Set-Row 18 $null $null "0.000008884"   "  +2.04%  " $true
Set-Row 19 $null $null "1.008"         "  -0.17%  " $true
Set-Row 20 $null $null $null           "  +0.48%  " $false
Set-Row 21 $null $null "27.072.68"     "  -2.42%  " $false
Set-Row 22 $null $null "5.167"         "  -1.38%  " $true
Set-Row 23 $null $null "10.59"         "  -0.35%  " $true
Set-Row 24 $null $null "2.099.01"      "  +0.01%  " $false
Set-Row 25 $null $null "152.90"        "  +0.82%  " $true
Set-Row 26 $null $null "1.845"         "  -1.93%  " $true
Set-Row 27 $null $null "18.28"         "  -1.20%  " $true
Set-Row 28 $null $null "2.094"         "  -1.74%  " $true
Set-Row 29 $null $null "5.130"         "  -0.95%  " $true
Set-Row 30 $null $null "115.71"        "  -0.62%  " $true
Set-Row 31 $null $null "0.08874"       "  -0.72%  " $true
Set-Row 32 $null $null "2.976"         "  +0.94%  " $true
Set-Row 33 $null $null "0.7301"        "  -1.42%  " $true
Set-Row 34 $null $null "4.447"         "  -1.37%  " $true
Set-Row 35 $null $null "1.140"         "  -2.01%  " $true
Set-Row 36 $null $null "2.502"         "  +4.00%  " $true

# Row 37 and 38 swap: VeChain/TrustWalletToken -> TrustWalletToken/VeChain
Set-Row 37 "TrustWalletToken" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt" "1.078" "  -0.95%  " $true
Set-Row 38 "VeChain" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet" "0.01954" "  +0.63%  " $true

Set-Row 39 $null $null "0.05235" "  -1.09%  " $true
Set-Row 40 $null $null $null     "  +0.67%  " $false
Set-Row 41 $null $null $null     "  -0.51%  " $false
Set-Row 42 $null $null "0.5205" "  -0.68%  " $true

# Row 43 and 44 swap: Frax/Algorand -> Algorand/Frax
Set-Row 43 "Algorand" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo" "0.1629" "  -0.83%  " $true
Set-Row 44 "Frax" "https://coinranking.com/coin/KfWtaeV1W+frax-frax" "0.8609" "  -14.72%  " $true

Set-Row 45 $null $null "8.246"  "  -1.73%  " $true
Set-Row 46 $null $null "0.4858" "  -0.36%  " $true

# Row 47 and 48 swap: EnergySwap/PaxDollar -> PaxDollar/EnergySwap
Set-Row 47 "PaxDollar" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp" "1.009" "  -0.04%  " $true
Set-Row 48 "EnergySwap" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens" "10.20" "  -1.46%  " $true

Set-Row 49 $null $null "102.87"  "  -1.52%  " $true
Set-Row 50 $null $null "1.628"   "  -1.26%  " $true
Set-Row 51 $null $null "0.06206" "  -1.37%  " $true
